# Applies the updated cryptos price/volume snapshot to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values that look numeric (e.g. '311.80') are written with a leading
# apostrophe so Excel stores them as text (quote-prefixed), matching the
# original inline-string cell type instead of converting them to numbers.

$ws.Range('D2').Value = '43.426.00'
$ws.Range('E2').Value = '  +1.71%  '
$ws.Range('D3').Value = '2.327.62'
$ws.Range('E3').Value = '  +1.06%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '''311.80'
$ws.Range('E5').Value = '  -1.39%  '
$ws.Range('D6').Value = '''108.26'
$ws.Range('E6').Value = '  +3.84%  '
$ws.Range('E7').Value = '  +0.83%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').Value = '''0.613'
$ws.Range('E9').Value = '  +1.77%  '
$ws.Range('D10').Value = '''40.67'
$ws.Range('E10').Value = '  +3.04%  '
$ws.Range('E11').Value = '  +1.15%  '
$ws.Range('D12').Value = '''8.57'
$ws.Range('E12').Value = '  +0.98%  '
$ws.Range('E13').Value = '  -1.22%  '
$ws.Range('E14').Value = '  -0.21%  '
$ws.Range('D15').Value = '''15.42'
$ws.Range('E15').Value = '  +0.79%  '
$ws.Range('D16').Value = '2.681.29'
$ws.Range('E16').Value = '  +1.09%  '
$ws.Range('D17').Value = '2.322.46'
$ws.Range('E17').Value = '  +0.87%  '
$ws.Range('D18').Value = '43.158.54'
$ws.Range('E18').Value = '  +0.77%  '
$ws.Range('E19').Value = '  +0.84%  '
$ws.Range('E20').Value = '  +0.98%  '
$ws.Range('D21').Value = '''12.97'
$ws.Range('E21').Value = '  -8.32%  '
$ws.Range('D22').Value = '''74.08'
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('E23').Value = '  -1.26%  '
$ws.Range('D24').Value = '''268.90'
$ws.Range('E24').Value = '  +2.03%  '
$ws.Range('E25').Value = '  +2.11%  '
$ws.Range('E26').Value = '  -0.39%  '
$ws.Range('D27').Value = '''7.62'
$ws.Range('E27').Value = '  +8.84%  '
$ws.Range('D28').Value = '''11.17'
$ws.Range('E28').Value = '  +2.55%  '
$ws.Range('D29').Value = '''2.29'
$ws.Range('E29').Value = '  -2.17%  '
$ws.Range('D30').Value = '''38.81'
$ws.Range('E30').Value = '  +3.78%  '
$ws.Range('D31').Value = '''22.61'
$ws.Range('E31').Value = '  +1.09%  '
$ws.Range('D32').Value = '''167.11'
$ws.Range('E32').Value = '  +0.30%  '
$ws.Range('D33').Value = '''0.0887'
$ws.Range('E33').Value = '  +1.59%  '
$ws.Range('E34').Value = '  +7.24%  '
$ws.Range('E35').Value = '  +0.52%  '
$ws.Range('E36').Value = '  +3.48%  '
$ws.Range('E37').Value = '  -1.54%  '
$ws.Range('D38').Value = '''0.0363'
$ws.Range('E38').Value = '  +3.89%  '
$ws.Range('E39').Value = '  +1.53%  '
$ws.Range('E40').Value = '  +5.49%  '
$ws.Range('E41').Value = '  +7.53%  '
$ws.Range('D42').Value = '''105.21'
$ws.Range('E42').Value = '  +13.60%  '
$ws.Range('D43').Value = '''71.64'
$ws.Range('E43').Value = '  +3.54%  '
$ws.Range('D44').Value = '''0.236'
$ws.Range('E44').Value = '  +2.92%  '
$ws.Range('D45').Value = '''13.23'
$ws.Range('E45').Value = '  +7.79%  '
$ws.Range('E46').Value = '  -0.36%  '
$ws.Range('D47').Value = '''114.06'
$ws.Range('E47').Value = '  -0.39%  '
$ws.Range('D48').Value = '1.660.36'
$ws.Range('E48').Value = '  -4.91%  '
$ws.Range('D49').Value = '''76.38'
$ws.Range('E49').Value = '  -5.26%  '
$ws.Range('E50').Value = '  +4.84%  '
$ws.Range('E51').Value = '  +1.51%  '
